$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1068.3334
$ws.Range("I28").Value = 1346.6666
$ws.Range("J28").Value = 882.7778
$ws.Range("K28").Value = 1346.6666
$ws.Range("L28").Value = 882.7778
$ws.Range("M28").Value = -861.6666
$ws.Range("N28").Value = -1852.7778
$ws.Range("H43").Value = 921.34485
$ws.Range("I43").Value = 366
$ws.Range("J43").Value = 1132.9048
$ws.Range("K43").Value = 366
$ws.Range("L43").Value = 1132.9048
$ws.Range("M43").Value = -297
$ws.Range("N43").Value = -1270.9048
$ws.Range("H112").Value = 22223794
$ws.Range("I112").Value = 200000320
$ws.Range("J112").Value = 1728.525
$ws.Range("K112").Value = 600000960
$ws.Range("L112").Value = 5185.575000000001
$ws.Range("M112").Value = -599999852
$ws.Range("N112").Value = -7401.575000000001
$ws.Range("H140").Value = 64998.57
$ws.Range("J140").Value = 64998.57
$ws.Range("L140").Value = 64998.57
$ws.Range("N140").Value = -75358.57000000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 45959.285
$ws.Range("J139").Value = 45959.285
$ws.Range("L139").Value = 45959.285
$ws.Range("N139").Value = -56239.285

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27496
$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -87480

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 37039428
$ws.Range("I16").Value = 1833.6666
$ws.Range("J16").Value = 55558228
$ws.Range("K16").Value = 1833.6666
$ws.Range("L16").Value = 55558228
$ws.Range("M16").Value = -1546.6666
$ws.Range("N16").Value = -55558802
$ws.Range("H31").Value = 2010.25
$ws.Range("I31").Value = 1407.0476
$ws.Range("K31").Value = 1407.0476
$ws.Range("M31").Value = -1112.0476
$ws.Range("H34").Value = 2010.25
$ws.Range("I34").Value = 1407.0476
$ws.Range("K34").Value = 1407.0476
$ws.Range("M34").Value = -1205.0476
$ws.Range("H96").Value = 25000
$ws.Range("J96").Value = 25000
$ws.Range("L96").Value = 25000
$ws.Range("N96").Value = -30492
$ws.Range("H107").Value = 1862.2
$ws.Range("I107").Value = 655.5
$ws.Range("J107").Value = 2666.6667
$ws.Range("K107").Value = 655.5
$ws.Range("L107").Value = 2666.6667
$ws.Range("M107").Value = 1264.5
$ws.Range("N107").Value = -6506.6667
$ws.Range("H113").Value = 37039428
$ws.Range("I113").Value = 1833.6666
$ws.Range("J113").Value = 55558228
$ws.Range("K113").Value = 1833.6666
$ws.Range("L113").Value = 55558228
$ws.Range("M113").Value = 336.3334
$ws.Range("N113").Value = -55562568
$ws.Range("H134").Value = 34457.574
$ws.Range("I134").Value = 2121.6843
$ws.Range("K134").Value = 6365.0529
$ws.Range("M134").Value = -3830.0529

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 934.931
$ws.Range("I107").Value = 727.5333000000001
$ws.Range("J107").Value = 1157.1428
$ws.Range("K107").Value = 2182.5999
$ws.Range("L107").Value = 3471.4284
$ws.Range("M107").Value = -262.5999000000002
$ws.Range("N107").Value = -7311.428400000001
$ws.Range("H110").Value = 8483.333000000001
$ws.Range("I110").Value = 900
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 2700
$ws.Range("L110").Value = 30000
$ws.Range("M110").Value = 1390
$ws.Range("N110").Value = -38180
$ws.Range("H115").Value = 2939.75
$ws.Range("I115").Value = 2603.6
$ws.Range("J115").Value = 3500
$ws.Range("K115").Value = 7810.799999999999
$ws.Range("L115").Value = 10500
$ws.Range("M115").Value = -6635.799999999999
$ws.Range("N115").Value = -12850

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 29630.5
$ws.Range("J39").Value = 29630.5
$ws.Range("L39").Value = 29630.5
$ws.Range("N39").Value = -30694.5
$ws.Range("H80").Value = 13692.4
$ws.Range("I80").Value = 51750
$ws.Range("J80").Value = 4178
$ws.Range("K80").Value = 51750
$ws.Range("L80").Value = 4178
$ws.Range("M80").Value = -50752
$ws.Range("N80").Value = -6174
$ws.Range("H83").Value = 13692.4
$ws.Range("I83").Value = 51750
$ws.Range("J83").Value = 4178
$ws.Range("K83").Value = 258750
$ws.Range("L83").Value = 20890
$ws.Range("M83").Value = -253758
$ws.Range("N83").Value = -30874
$ws.Range("H107").Value = 1625
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 1666.6666
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 1666.6666
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -5506.6666
$ws.Range("H113").Value = 1715.15
$ws.Range("I113").Value = 1664.8182
$ws.Range("J113").Value = 1776.6666
$ws.Range("K113").Value = 1664.8182
$ws.Range("L113").Value = 1776.6666
$ws.Range("M113").Value = 505.1818000000001
$ws.Range("N113").Value = -6116.6666
$ws.Range("H114").Value = 29797.059
$ws.Range("I114").Value = 15000
$ws.Range("J114").Value = 30721.875
$ws.Range("K114").Value = 15000
$ws.Range("L114").Value = 30721.875
$ws.Range("M114").Value = -10661
$ws.Range("N114").Value = -39399.875
$ws.Range("H132").Value = 202077.7
$ws.Range("I132").Value = 167799.83
$ws.Range("K132").Value = 503399.49
$ws.Range("M132").Value = -500869.49

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 744.6923
$ws.Range("I46").Value = 723.4167
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 723.4167
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -535.4167
$ws.Range("N46").Value = -1376
$ws.Range("H61").Value = 3075.25
$ws.Range("I61").Value = 3328.5715
$ws.Range("J61").Value = 1302
$ws.Range("K61").Value = 3328.5715
$ws.Range("L61").Value = 1302
$ws.Range("M61").Value = -3126.5715
$ws.Range("N61").Value = -1706
$ws.Range("H113").Value = 3075.25
$ws.Range("I113").Value = 3328.5715
$ws.Range("J113").Value = 1302
$ws.Range("K113").Value = 3328.5715
$ws.Range("L113").Value = 1302
$ws.Range("M113").Value = -1158.5715
$ws.Range("N113").Value = -5642
$ws.Range("H122").Value = 4031.074
$ws.Range("I122").Value = 3686.4167
$ws.Range("J122").Value = 4306.8
$ws.Range("K122").Value = 11059.2501
$ws.Range("L122").Value = 12920.4
$ws.Range("M122").Value = -8609.250100000001
$ws.Range("N122").Value = -17820.4
$ws.Range("H136").Value = 223663.56
$ws.Range("I136").Value = 126621.5
$ws.Range("K136").Value = 379864.5
$ws.Range("M136").Value = -377314.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 220.71428
$ws.Range("I107").Value = 199.16667
$ws.Range("K107").Value = 597.50001
$ws.Range("M107").Value = 1322.49999
$ws.Range("H113").Value = 1656.25
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1656.25
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4968.75
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9308.75
